# Update "想去人数" (column F) figures across the 展览 / 演出 / 本地生活 / 全部类型
# sheets to match the freshly scraped counts (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 151
$ws.Cells.Item(3, 6).Value = 198
$ws.Cells.Item(6, 6).Value = 1315
$ws.Cells.Item(7, 6).Value = 70
$ws.Cells.Item(9, 6).Value = 390
$ws.Cells.Item(10, 6).Value = 446
$ws.Cells.Item(11, 6).Value = 817
$ws.Cells.Item(12, 6).Value = 211
$ws.Cells.Item(13, 6).Value = 749
$ws.Cells.Item(14, 6).Value = 315
$ws.Cells.Item(15, 6).Value = 469
$ws.Cells.Item(17, 6).Value = 1050
$ws.Cells.Item(18, 6).Value = 494
$ws.Cells.Item(20, 6).Value = 410
$ws.Cells.Item(21, 6).Value = 103
$ws.Cells.Item(22, 6).Value = 221
$ws.Cells.Item(24, 6).Value = 56
$ws.Cells.Item(25, 6).Value = 490
$ws.Cells.Item(26, 6).Value = 441
$ws.Cells.Item(27, 6).Value = 290

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 379
$ws.Cells.Item(5, 6).Value = 49
$ws.Cells.Item(6, 6).Value = 47
$ws.Cells.Item(7, 6).Value = 293
$ws.Cells.Item(11, 6).Value = 160
$ws.Cells.Item(12, 6).Value = 145

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 355

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 355
$ws.Cells.Item(4, 6).Value = 151
$ws.Cells.Item(5, 6).Value = 198
$ws.Cells.Item(8, 6).Value = 1315
$ws.Cells.Item(10, 6).Value = 70
$ws.Cells.Item(11, 6).Value = 379
$ws.Cells.Item(13, 6).Value = 49
$ws.Cells.Item(14, 6).Value = 390
$ws.Cells.Item(15, 6).Value = 47
$ws.Cells.Item(16, 6).Value = 293
$ws.Cells.Item(17, 6).Value = 446
$ws.Cells.Item(18, 6).Value = 817
$ws.Cells.Item(19, 6).Value = 211
$ws.Cells.Item(20, 6).Value = 749
$ws.Cells.Item(21, 6).Value = 315
$ws.Cells.Item(22, 6).Value = 469
$ws.Cells.Item(24, 6).Value = 1050
$ws.Cells.Item(25, 6).Value = 494
$ws.Cells.Item(29, 6).Value = 410
$ws.Cells.Item(31, 6).Value = 103
$ws.Cells.Item(32, 6).Value = 160
$ws.Cells.Item(33, 6).Value = 221
$ws.Cells.Item(35, 6).Value = 56
$ws.Cells.Item(36, 6).Value = 145
$ws.Cells.Item(38, 6).Value = 490
$ws.Cells.Item(41, 6).Value = 441
$ws.Cells.Item(42, 6).Value = 290
